$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing data rows 2-21 down to 3-22),
# then strip the formatting Excel's row-insert copies down from row 1 so the
# new row matches the plain (unstyled) numeric cells used elsewhere in the sheet.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:C2").ClearFormats()

# Populate the newly inserted row 2 with its values.
$ws.Range("A2").Value = -0.0303905457258224
$ws.Range("B2").Value = -0.02122756652534
$ws.Range("C2").Value = -0.038026362657547

# Append the new trailing rows (23-31) with their values.
$newRows = @(
    @(-0.0070249503478407, 4.270253658294678, -0.0296269636601209),
    @(0.4137084782123565, 2.936276435852051, 0.2823724448680877),
    @(0.0429132841527462, 1.122159481048584, 0.1867720484733581),
    @(0.06475171446800231, -1.842216849327088, -0.6108652353286743),
    @(0.0862847194075584, -5.713422775268555, -1.346194267272949),
    @(-0.1818851232528686, -4.851491928100586, 1.392772793769836),
    @(-0.3181080818176269, -3.869678497314453, 0.9886853694915771),
    @(0.1050688251852989, -2.216677188873291, 0.3729332387447357),
    @(0.1996002197265625, 1.434922456741333, -0.2237294018268585)
)

$r = 23
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $r = $r + 1
}
